$d = $word.ActiveDocument

function FindParaIndex($text) {
  for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    $t = $t.TrimEnd([char]13)
    if ($t -eq $text) {
      return $i
    }
  }
  return -1
}

function InsertRoleAfter($anchorText, $roleText) {
  $idx = FindParaIndex $anchorText
  $d.Paragraphs.Item($idx).Range.InsertParagraphAfter()
  $d.Paragraphs.Item($idx + 1).Range.Text = $roleText
}

# Insert the new "Role:" lines after each team line, top to bottom.
# (FindParaIndex re-scans each time, so earlier insertions don't
#  invalidate the indices used for later ones.)
InsertRoleAfter "Team Leader: K.Keerthana&kumarkeerthana2692@gmail.com" "Role:Project document "
InsertRoleAfter "Keerthana.S&keerthana2007ss@gmail.com" "Role:Coding creation"
InsertRoleAfter "Keerthana.S&keerkeerthana2996@gmail.com " "Role:Coding creation"
InsertRoleAfter "Kokila.P&kokilapalani84@gmail.com" "Role:Demo Link"
InsertRoleAfter "Komathi.S&kkomath517@gamil.com" "Role:Project document "

# Relocate the _GoBack bookmark from the end of the "Team Id" paragraph to
# a collapsed point inside the new "Role:Coding creation" paragraph that
# follows "Keerthana.S&keerthana2007ss@gmail.com" (splitting its text into
# "Role:C" | bookmark | "oding creation").
$d.Bookmarks.Item("_GoBack").Delete()
$roleIdx = FindParaIndex "Role:Coding creation"
$rolePara = $d.Paragraphs.Item($roleIdx)
$splitPoint = $rolePara.Range.Start + 6
$target = $d.Range($splitPoint, $splitPoint)
$d.Bookmarks.Add("_GoBack", $target)
